# Apply "Avrg" header column, Micro row (54) literal values, and
# "Macro wo reviews" row (56) formulas to every results sheet in the
# workbook, and tidy up the view/selection + row heights to match.

$wb = $excel.ActiveWorkbook

# Per-sheet Micro-average literals (B54, C54, D54) taken from the commit.
$microValues = @{
    "strict_strict"   = @(0.695, 0.796, 0.66)
    "strict_relaxed"  = @(0.83,  0.874, 0.824)
    "partial_strict"  = @(0.738, 0.842, 0.701)
    "partial_relaxed" = @(0.904, 0.939, 0.897)
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $name = $ws.Name

    # --- Header: new "Avrg" column label in E1 ---------------------------
    $ws.Range("E1").Value2 = "Avrg"

    # --- Row 53 ("Macro"): give E53 the same 0.000 number format as the
    #     other averages in that row (previously General) -----------------
    $ws.Range("E53").NumberFormat = "0.000"

    # --- Row 54 ("Micro"): new literal percentages + average formula -----
    $vals = $microValues[$name]
    $ws.Range("B54").Value2 = $vals[0]
    $ws.Range("C54").Value2 = $vals[1]
    $ws.Range("D54").Value2 = $vals[2]
    $ws.Range("E54").Formula = "=SUM(B54:D54)/3"
    if ($name -eq "strict_strict") {
        $ws.Range("E54").NumberFormat = "General"
    } else {
        $ws.Range("E54").NumberFormat = "0.000"
    }

    # --- Row 56 ("Macro wo reviews"): averages excluding the two review
    #     papers (rows 39 and 42) -----------------------------------------
    $ws.Range("B56").Formula = "=(SUM(B2:B51) - B39 - B42 )/48"
    $ws.Range("C56").Formula = "=(SUM(C2:C51) - C39 - C42 )/48"
    $ws.Range("D56").Formula = "=(SUM(D2:D51) - D39 - D42 )/48"
    $ws.Range("E56").Formula = "=SUM(B56:D56)/3"
    $ws.Range("B56:E56").NumberFormat = "0.000"

    # --- Row heights: header, Macro, Micro and Macro-wo-reviews rows grow
    #     slightly (12.8 -> 13.85) to fit the updated layout ---------------
    $ws.Rows.Item(1).RowHeight = 13.85
    $ws.Rows.Item(53).RowHeight = 13.85
    $ws.Rows.Item(54).RowHeight = 13.85
    $ws.Rows.Item(56).RowHeight = 13.85

    # --- View: scroll back to the top and select E1 instead of E2:E51 ----
    $ws.Activate()
    $excel.ActiveWindow.ScrollRow = 1
    $excel.ActiveWindow.ScrollColumn = 1
    $ws.Range("E1").Select()
}

# Restore the originally active tab (partial_relaxed, the 4th sheet).
$wb.Worksheets.Item(4).Activate()
$wb.Worksheets.Item(4).Range("E1").Select()
